$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.205.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").Value = "'1.844.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'240.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").Value = "'0.6717"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.87%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'0.07424"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("E9").Value = "  -1.96%  "

$ws.Range("D10").Value = "'22.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("D11").Value = "'0.07717"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").Value = "'1.815.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("D13").Value = "'5.006"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("D14").Value = "'0.6746"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("D15").Value = "'86.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("D16").Value = "'6.129"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").Value = "'29.168.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").Value = "'0.000008307"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.65%  "

$ws.Range("D19").Value = "'228.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").Value = "'7.187"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.21%  "

$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").Value = "'160.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("E25").Value = "  -0.66%  "

$ws.Range("D26").Value = "'0.1405"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.71%  "

$ws.Range("D27").Value = "'18.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").Value = "'1.510"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("D29").Value = "'4.176"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.59%  "

$ws.Range("D30").Value = "'4.069"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.01%  "

$ws.Range("D31").Value = "'1.188"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.74%  "

$ws.Range("D32").Value = "'0.05318"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.89%  "

$ws.Range("D33").Value = "'0.7582"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").Value = "'1.875"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.58%  "

$ws.Range("E35").Value = "  +0.44%  "

$ws.Range("D36").Value = "'2.679"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").Value = "'1.331.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.84%  "

$ws.Range("D38").Value = "'0.01800"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.05%  "

$ws.Range("D39").Value = "'2.731"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("D40").Value = "'0.9242"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.78%  "

$ws.Range("D41").Value = "'5.958"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.58%  "

$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "'103.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("D44").Value = "'0.07921"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.60%  "

$ws.Range("D45").Value = "'1.969.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("D46").Value = "'0.5162"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "

$ws.Range("D47").Value = "'1.772"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("E48").Value = "  -2.18%  "

$ws.Range("D49").Value = "'63.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.25%  "

$ws.Range("D50").Value = "'9.146"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.35%  "

$ws.Range("D51").Value = "'0.05947"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
